$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.202.15"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.635.61"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.97"
$ws.Range("E5").Value = "  +4.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.64"
$ws.Range("E6").Value = "  +1.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.70"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  +6.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.350"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.094.86"
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "61.176.44"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.88"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000146"
$ws.Range("E16").Value = "  +3.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.633.68"
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.80"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "357.32"
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.69"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.26"
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.83"
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.432"
$ws.Range("E24").Value = "  +2.87%  "
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.743.21"
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0872"
$ws.Range("E28").Value = "  +3.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.44"
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  +8.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.61"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.63"
$ws.Range("E33").Value = "  +4.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "151.50"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.20"
$ws.Range("E35").Value = "  +4.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.21"
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.922"
$ws.Range("E37").Value = "  +9.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.888"
$ws.Range("E38").Value = "  +3.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.51"
$ws.Range("E39").Value = "  +2.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.83"
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "294.48"
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.644"
$ws.Range("E42").Value = "  +4.22%  "
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0564"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.10"
$ws.Range("E45").Value = "  +5.97%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.87"
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("E48").Value = "  +3.20%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.88"
$ws.Range("E49").Value = "  +5.83%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.34"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.10"
$ws.Range("E51").Value = "  +6.76%  "
